# Updated cryptos list on Tue Oct  3 02:41:02 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# table with new scraped figures, and fixes the rank ordering of two
# coin pairs that swapped places (WEMIXToken/PaxDollar @ rows 41-42,
# RocketPoolETH/MXToken @ rows 44-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a price into column D while preserving it as plain text.
# Several prices (e.g. "215.29", "1.02") are syntactically valid numbers,
# and a bare Range.Value assignment would let Excel auto-convert them to
# numerics (dropping the text formatting the sheet relies on). Briefly
# forcing a Text number format keeps the write as a string, then the
# style is reset back to Normal so no stray formatting is left behind.
function Set-PriceText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# row 2 - Bitcoin
$ws.Range("D2").Value = "27.562.26"
$ws.Range("E2").Value = "  -1.45%  "

# row 3 - Ethereum
$ws.Range("D3").Value = "1.666.81"
$ws.Range("E3").Value = "  -3.38%  "

# row 5 - BNB
Set-PriceText "D5" "215.29"
$ws.Range("E5").Value = "  -1.54%  "

# row 6 - XRP
$ws.Range("E6").Value = "  -1.68%  "

# row 8 - Solana
Set-PriceText "D8" "23.56"
$ws.Range("E8").Value = "  -1.88%  "

# row 9 - Cardano
$ws.Range("E9").Value = "  -1.16%  "

# row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.80%  "

# row 11 - TRON
$ws.Range("E11").Value = "  -2.54%  "

# row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.902.46"
$ws.Range("E12").Value = "  -3.38%  "

# row 13 - WrappedEther
$ws.Range("D13").Value = "1.664.88"
$ws.Range("E13").Value = "  -3.34%  "

# row 14 - Polkadot
$ws.Range("E14").Value = "  -2.54%  "

# row 15 - Polygon
$ws.Range("E15").Value = "  -2.15%  "

# row 16 - Litecoin
Set-PriceText "D16" "66.25"
$ws.Range("E16").Value = "  -2.46%  "

# row 17 - BitcoinCash
Set-PriceText "D17" "251.18"
$ws.Range("E17").Value = "  +2.70%  "

# row 18 - WrappedBTC
$ws.Range("D18").Value = "27.575.44"
$ws.Range("E18").Value = "  -1.18%  "

# row 19 - ShibaInu
$ws.Range("E19").Value = "  -3.07%  "

# row 20 - Chainlink
Set-PriceText "D20" "7.53"
$ws.Range("E20").Value = "  -4.57%  "

# row 21 - Dai
$ws.Range("E21").Value = "  -0.03%  "

# row 23 - Avalanche
Set-PriceText "D23" "9.30"
$ws.Range("E23").Value = "  -4.56%  "

# row 25 - Monero
Set-PriceText "D25" "146.57"
$ws.Range("E25").Value = "  -1.95%  "

# row 26 - EthereumClassic
Set-PriceText "D26" "16.57"
$ws.Range("E26").Value = "  -1.44%  "

# row 27 - Cosmos
$ws.Range("E27").Value = "  -4.89%  "

# row 28 - Stellar
$ws.Range("E28").Value = "  -2.21%  "

# row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.08%  "

# row 30 - PancakeSwap
$ws.Range("E30").Value = "  +3.96%  "

# row 31 - Hedera
$ws.Range("E31").Value = "  -0.49%  "

# row 33 - Maker
$ws.Range("D33").Value = "1.477.38"
$ws.Range("E33").Value = "  -0.99%  "

# row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -5.51%  "

# row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -5.38%  "

# row 36 - ARBITRUM
$ws.Range("E36").Value = "  -1.91%  "

# row 37 - HuobiToken
$ws.Range("E37").Value = "  -0.92%  "

# row 38 - ImmutableX
Set-PriceText "D38" "0.578"
$ws.Range("E38").Value = "  -5.80%  "

# row 39 - VeChain
$ws.Range("E39").Value = "  -2.32%  "

# row 40 - Aave
$ws.Range("E40").Value = "  -2.36%  "

# row 41 - was PaxDollar, now WEMIXToken (swapped with row 42)
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-PriceText "D41" "1.02"
$ws.Range("E41").Value = "  -4.58%  "

# row 42 - was WEMIXToken, now PaxDollar (swapped with row 41)
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-PriceText "D42" "1.00"
$ws.Range("E42").Value = "  +0.05%  "

# row 43 - FraxShare
$ws.Range("E43").Value = "  -6.99%  "

# row 44 - was MXToken, now RocketPoolETH (swapped with row 45)
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.810.54"
$ws.Range("E44").Value = "  -3.34%  "

# row 45 - was RocketPoolETH, now MXToken (swapped with row 44)
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-PriceText "D45" "2.21"
$ws.Range("E45").Value = "  -3.30%  "

# row 46 - TrustWalletToken
Set-PriceText "D46" "0.792"
$ws.Range("E46").Value = "  +0.03%  "

# row 47 - RenderToken
$ws.Range("E47").Value = "  -2.02%  "

# row 48 - Quant
Set-PriceText "D48" "89.48"
$ws.Range("E48").Value = "  -1.83%  "

# row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -1.60%  "

# row 50 - BitcoinSV
Set-PriceText "D50" "42.15"
$ws.Range("E50").Value = "  +15.77%  "

# row 51 - Algorand
$ws.Range("E51").Value = "  -3.34%  "
